$wb = $excel.ActiveWorkbook

# --- Sheet "task 3" (internal sheet1.xml) holds the K-column concatenation formulas ---
$ws3 = $wb.Worksheets.Item("task 3")

# New formula: the trailing two SUBSTITUTE() calls that used to replace the
# space between date and time with an underscore now replace it with a
# plain space (i.e. keep it as-is). NOTE: deliberately built with plain
# string concatenation (+) instead of the -f operator, since -f mangles
# the embedded Cyrillic characters in this runtime.
for ($r = 1; $r -le 20; $r++) {
    $formula = '=SUBSTITUTE(A' + $r + '," ","_")&" "&SUBSTITUTE(B' + $r + '," ","_")&" "&C' + $r + '&" "&D' + $r + '&" "&E' + $r + '&" "&F' + $r + '&" "&SUBSTITUTE(TEXT(G' + $r + ',"0,00"),",",".")&" "&H' + $r + '&" "&SUBSTITUTE(TEXT(I' + $r + ',"ДД:ММ:ГГГГ_ЧЧ:ММ:СС")," "," ")&" "&SUBSTITUTE(TEXT(J' + $r + ',"ДД:ММ:ГГГГ_ЧЧ:ММ:СС")," "," ")'
    $ws3.Range("K$r").Formula = $formula
}

# --- Sheet "task 5" (internal sheet3.xml) loses the active-tab flag and gets
#     a new selection (activeCell G15) ---
$ws5 = $wb.Worksheets.Item("task 5")
$ws5.Range("G15").Select() | Out-Null

# --- Back on "task 3": it becomes the active tab, with K1:K20 selected
#     (active cell K1). Doing this last makes "task 3" the workbook's
#     active/visible sheet when saved, matching the target bookViews. ---
$ws3.Activate()
$ws3.Range("K1:K20").Select() | Out-Null
